$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit3"
$ws.Range("C2").Value = "Robo1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.044914666666666
$ws.Range("H2").Value = 3.134744
$ws.Range("I2").Value = 0.006668841574421894
$ws.Range("J2").Value = 0.006668841574421893
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1757713333333334
$ws.Range("N2").Value = 0.5273140000000001
$ws.Range("O2").Value = 0.009444264308298454
$ws.Range("P2").Value = 0.009444264308298454
$ws.Range("Q2").Value = 0.1836660441795555
$ws.Range("R2").Value = 1.652994397616
$ws.Range("S2").Value = 0.00006298230245900955344764971
$ws.Range("T2").Value = 0.00006298230245900955344764971

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit3"
$ws.Range("C3").Value = "Robo1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.044914666666666
$ws.Range("H3").Value = 3.134744
$ws.Range("I3").Value = 0.006668841574421894
$ws.Range("J3").Value = 0.006668841574421893
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.98312366666667
$ws.Range("N3").Value = 50.949371
$ws.Range("O3").Value = 0.9125100529581165
$ws.Range("P3").Value = 0.9125100529581165
$ws.Range("Q3").Value = 17.74591500511378
$ws.Range("R3").Value = 159.713235046024
$ws.Range("S3").Value = 0.006085384978245011
$ws.Range("T3").Value = 0.00608538497824501

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit3"
$ws.Range("C4").Value = "Robo1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.044914666666666
$ws.Range("H4").Value = 3.134744
$ws.Range("I4").Value = 0.006668841574421894
$ws.Range("J4").Value = 0.006668841574421893
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.452542333333333
$ws.Range("N4").Value = 4.357627
$ws.Range("O4").Value = 0.07804568273358503
$ws.Range("P4").Value = 0.07804568273358505
$ws.Range("Q4").Value = 1.517782788054222
$ws.Range("R4").Value = 13.660045092488
$ws.Range("S4").Value = 0.0005204742937178728
$ws.Range("T4").Value = 0.0005204742937178729

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit3"
$ws.Range("C5").Value = "Robo1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 127.867017
$ws.Range("H5").Value = 383.601051
$ws.Range("I5").Value = 0.8160713081836135
$ws.Range("J5").Value = 0.8160713081836134
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1757713333333334
$ws.Range("N5").Value = 0.5273140000000001
$ws.Range("O5").Value = 0.009444264308298454
$ws.Range("P5").Value = 0.009444264308298454
$ws.Range("Q5").Value = 22.475356067446
$ws.Range("R5").Value = 202.278204607014
$ws.Range("S5").Value = 0.007707193128904929
$ws.Range("T5").Value = 0.007707193128904928

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit3"
$ws.Range("C6").Value = "Robo1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 127.867017
$ws.Range("H6").Value = 383.601051
$ws.Range("I6").Value = 0.8160713081836135
$ws.Range("J6").Value = 0.8160713081836134
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.98312366666667
$ws.Range("N6").Value = 50.949371
$ws.Range("O6").Value = 0.9125100529581165
$ws.Range("P6").Value = 0.9125100529581165
$ws.Range("Q6").Value = 2171.581362598769
$ws.Range("R6").Value = 19544.23226338892
$ws.Range("S6").Value = 0.7446732726482286
$ws.Range("T6").Value = 0.7446732726482285

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit3"
$ws.Range("C7").Value = "Robo1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 127.867017
$ws.Range("H7").Value = 383.601051
$ws.Range("I7").Value = 0.8160713081836135
$ws.Range("J7").Value = 0.8160713081836134
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.452542333333333
$ws.Range("N7").Value = 4.357627
$ws.Range("O7").Value = 0.07804568273358503
$ws.Range("P7").Value = 0.07804568273358505
$ws.Range("Q7").Value = 185.732255229553
$ws.Range("R7").Value = 1671.590297065977
$ws.Range("S7").Value = 0.06369084240648
$ws.Range("T7").Value = 0.06369084240648

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit3"
$ws.Range("C8").Value = "Robo1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 27.77415166666666
$ws.Range("H8").Value = 83.32245499999999
$ws.Range("I8").Value = 0.1772598502419647
$ws.Range("J8").Value = 0.1772598502419647
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1757713333333334
$ws.Range("N8").Value = 0.5273140000000001
$ws.Range("O8").Value = 0.009444264308298454
$ws.Range("P8").Value = 0.009444264308298454
$ws.Range("Q8").Value = 4.881899670652222
$ws.Range("R8").Value = 43.93709703587
$ws.Range("S8").Value = 0.001674088876934516
$ws.Range("T8").Value = 0.001674088876934516

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit3"
$ws.Range("C9").Value = "Robo1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 27.77415166666666
$ws.Range("H9").Value = 83.32245499999999
$ws.Range("I9").Value = 0.1772598502419647
$ws.Range("J9").Value = 0.1772598502419647
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.98312366666667
$ws.Range("N9").Value = 50.949371
$ws.Range("O9").Value = 0.9125100529581165
$ws.Range("P9").Value = 0.9125100529581165
$ws.Range("Q9").Value = 471.6918524917561
$ws.Range("R9").Value = 4245.226672425804
$ws.Range("S9").Value = 0.161751395331643
$ws.Range("T9").Value = 0.161751395331643

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit3"
$ws.Range("C10").Value = "Robo1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 27.77415166666666
$ws.Range("H10").Value = 83.32245499999999
$ws.Range("I10").Value = 0.1772598502419647
$ws.Range("J10").Value = 0.1772598502419647
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.452542333333333
$ws.Range("N10").Value = 4.357627
$ws.Range("O10").Value = 0.07804568273358503
$ws.Range("P10").Value = 0.07804568273358505
$ws.Range("Q10").Value = 40.34313106825389
$ws.Range("R10").Value = 363.088179614285
$ws.Range("S10").Value = 0.01383436603338717
$ws.Range("T10").Value = 0.01383436603338717
